# Updates cryptos list values (price + volume columns, plus a few
# Coin/Link reorderings) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new literal text value. Price
# values (column D) sometimes parse as plain numbers (e.g. "574.75"),
# so force those cells to Text format first -- exactly like the source
# data, which stores every price/volume/coin/link cell as text -- to
# avoid Excel auto-converting them to numeric cells on input.
$updates = @(
    @{ Cell = 'D2'; Value = '64.421.86' }
    @{ Cell = 'E2'; Value = '  -0.76%  ' }
    @{ Cell = 'D3'; Value = '3.447.71' }
    @{ Cell = 'E3'; Value = '  +0.44%  ' }
    @{ Cell = 'E4'; Value = '  -0.03%  ' }
    @{ Cell = 'D5'; Value = '574.75' }
    @{ Cell = 'E5'; Value = '  -0.04%  ' }
    @{ Cell = 'D6'; Value = '164.93' }
    @{ Cell = 'E6'; Value = '  +3.61%  ' }
    @{ Cell = 'E7'; Value = '  +0.01%  ' }
    @{ Cell = 'D8'; Value = '3.450.14' }
    @{ Cell = 'E8'; Value = '  +0.45%  ' }
    @{ Cell = 'D9'; Value = '0.558' }
    @{ Cell = 'E9'; Value = '  -4.91%  ' }
    @{ Cell = 'D10'; Value = '7.31' }
    @{ Cell = 'E10'; Value = '  +1.16%  ' }
    @{ Cell = 'D11'; Value = '0.121' }
    @{ Cell = 'E11'; Value = '  -1.08%  ' }
    @{ Cell = 'D12'; Value = '0.428' }
    @{ Cell = 'E12'; Value = '  -4.44%  ' }
    @{ Cell = 'D13'; Value = '4.038.07' }
    @{ Cell = 'E13'; Value = '  +0.32%  ' }
    @{ Cell = 'D14'; Value = '0.136' }
    @{ Cell = 'E14'; Value = '  +1.32%  ' }
    @{ Cell = 'D15'; Value = '27.46' }
    @{ Cell = 'E15'; Value = '  -1.14%  ' }
    @{ Cell = 'D16'; Value = '0.0000176' }
    @{ Cell = 'E16'; Value = '  -6.42%  ' }
    @{ Cell = 'D17'; Value = '64.457.72' }
    @{ Cell = 'E17'; Value = '  -0.76%  ' }
    @{ Cell = 'D18'; Value = '3.433.12' }
    @{ Cell = 'E18'; Value = '  +0.76%  ' }
    @{ Cell = 'D19'; Value = '6.20' }
    @{ Cell = 'E19'; Value = '  -2.91%  ' }
    @{ Cell = 'D20'; Value = '13.73' }
    @{ Cell = 'E20'; Value = '  -0.99%  ' }
    @{ Cell = 'D21'; Value = '379.99' }
    @{ Cell = 'E21'; Value = '  -0.43%  ' }
    @{ Cell = 'D22'; Value = '7.90' }
    @{ Cell = 'E22'; Value = '  -0.88%  ' }
    @{ Cell = 'D23'; Value = '0.999' }
    @{ Cell = 'E23'; Value = '  -0.10%  ' }
    @{ Cell = 'D24'; Value = '71.51' }
    @{ Cell = 'E24'; Value = '  -0.83%  ' }
    @{ Cell = 'D25'; Value = '0.522' }
    @{ Cell = 'E25'; Value = '  -4.79%  ' }
    @{ Cell = 'D26'; Value = '0.0000118' }
    @{ Cell = 'E26'; Value = '  -1.22%  ' }
    @{ Cell = 'D27'; Value = '9.64' }
    @{ Cell = 'E27'; Value = '  -3.14%  ' }
    @{ Cell = 'D28'; Value = '0.178' }
    @{ Cell = 'E28'; Value = '  +0.09%  ' }
    @{ Cell = 'D29'; Value = '1.00' }
    @{ Cell = 'E29'; Value = '  +0.12%  ' }
    @{ Cell = 'D30'; Value = '6.15' }
    @{ Cell = 'E30'; Value = '  -0.02%  ' }
    @{ Cell = 'D31'; Value = '1.42' }
    @{ Cell = 'E31'; Value = '  -2.89%  ' }
    @{ Cell = 'D33'; Value = '23.08' }
    @{ Cell = 'E33'; Value = '  -0.88%  ' }
    @{ Cell = 'D34'; Value = '7.21' }
    @{ Cell = 'E34'; Value = '  +2.35%  ' }
    @{ Cell = 'D35'; Value = '1.52' }
    @{ Cell = 'E35'; Value = '  -3.81%  ' }
    @{ Cell = 'D36'; Value = '160.25' }
    @{ Cell = 'E36'; Value = '  -0.31%  ' }
    @{ Cell = 'D37'; Value = '0.865' }
    @{ Cell = 'E37'; Value = '  +11.75%  ' }
    @{ Cell = 'D38'; Value = '1.84' }
    @{ Cell = 'E38'; Value = '  -3.37%  ' }
    @{ Cell = 'D39'; Value = '2.828.08' }
    @{ Cell = 'E39'; Value = '  -2.37%  ' }
    @{ Cell = 'D40'; Value = '0.0733' }
    @{ Cell = 'E40'; Value = '  -2.31%  ' }
    @{ Cell = 'D41'; Value = '26.15' }
    @{ Cell = 'E41'; Value = '  -0.81%  ' }
    @{ Cell = 'B42'; Value = 'RenderToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D42'; Value = '6.53' }
    @{ Cell = 'E42'; Value = '  -3.38%  ' }
    @{ Cell = 'B43'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D43'; Value = '26.55' }
    @{ Cell = 'E43'; Value = '  +2.43%  ' }
    @{ Cell = 'B44'; Value = 'OKB' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D44'; Value = '42.99' }
    @{ Cell = 'E44'; Value = '  -0.62%  ' }
    @{ Cell = 'D45'; Value = '4.46' }
    @{ Cell = 'E45'; Value = '  -2.22%  ' }
    @{ Cell = 'B46'; Value = 'dogwifhat' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' }
    @{ Cell = 'D46'; Value = '2.53' }
    @{ Cell = 'E46'; Value = '  +11.42%  ' }
    @{ Cell = 'B47'; Value = 'VeChain' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D47'; Value = '0.0309' }
    @{ Cell = 'E47'; Value = '  -2.55%  ' }
    @{ Cell = 'D48'; Value = '339.80' }
    @{ Cell = 'E48'; Value = '  +7.40%  ' }
    @{ Cell = 'E49'; Value = '  -1.65%  ' }
    @{ Cell = 'D50'; Value = '0.104' }
    @{ Cell = 'E50'; Value = '  -2.09%  ' }
    @{ Cell = 'D51'; Value = '6.36' }
    @{ Cell = 'E51'; Value = '  -2.53%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $text = $u.Value
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain number (single decimal point) -- Excel would
        # otherwise silently coerce it to a numeric cell on assignment, so
        # force Text format first to keep it a string, matching the source.
        $range.NumberFormat = '@'
    }
    $range.Value = $text
}
